$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the "Absent" column (H) for rows where it was missing/incorrect,
# completing the consolidated attendance report.
$ws.Range("H3").Value = 1
$ws.Range("H5").Value = 0
$ws.Range("H12").Value = 1
$ws.Range("H14").Value = 0
